$wb = $excel.ActiveWorkbook

# --- Sheet1 (Task13): add row 3 "Welcome to UiPath Training Session" ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("A3").Value = "Welcome "
$ws1.Range("B3").Value = "to "
$ws1.Range("C3").Value = "UiPath "
$ws1.Range("D3").Value = "Training "
$ws1.Range("E3").Value = "Session"

# --- Sheet2 (Task14): add headers Column1/2/3 in row 1 and row 3 text ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("C1").Value = "Column1"
$ws2.Range("D1").Value = "Column2"
$ws2.Range("E1").Value = "Column3"
$ws2.Range("A3").Value = "Welcome "
$ws2.Range("B3").Value = "to "
$ws2.Range("C3").Value = "UiPath "
$ws2.Range("D3").Value = "Training "
$ws2.Range("E3").Value = "Session"

# --- Sheet3: same as Sheet2 ---
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("C1").Value = "Column1"
$ws3.Range("D1").Value = "Column2"
$ws3.Range("E1").Value = "Column3"
$ws3.Range("A3").Value = "Welcome "
$ws3.Range("B3").Value = "to "
$ws3.Range("C3").Value = "UiPath "
$ws3.Range("D3").Value = "Training "
$ws3.Range("E3").Value = "Session"

# --- Restore final active sheet/selection: Sheet1 selection moves to C4,
#     Sheet2 stays the active/selected tab ---
$ws1.Range("C4").Select()
$ws2.Select()
